$wb = $excel.ActiveWorkbook

# --- 1. Update the "last updated" date (ERT_ATFM_YY!B2), the other sheets
#        (MM/FAB/LOC) reference this cell via formula so they refresh too.
$wsYY = $wb.Worksheets.Item("ERT_ATFM_YY")
$wsYY.Range("B2").Value = 45110

# --- 2. FABEC (row12) / SW FAB (row14) minutes of delay, ERT_ATFM_FAB sheet.
#        Column E on that sheet is a shared formula (=D/C) so it recalculates
#        automatically once D is edited.
$wsFAB = $wb.Worksheets.Item("ERT_ATFM_FAB")
$wsFAB.Range("D12").Value = 10466924
$wsFAB.Range("D14").Value = 1002313

# --- 3. DSNA (row13) / ENAIRE (row15) minutes of delay, ERT_ATFM_LOC sheet.
$wsLOC = $wb.Worksheets.Item("ERT_ATFM_LOC")
$wsLOC.Range("D13").Value = 4342838
$wsLOC.Range("D15").Value = 598117

Write-Host "numbers updated"

# --- 4. "Change Log" sheet: fill in the three new change-log rows describing
#        the ENAIRE -> DSNA delay reallocations.
$wsLog = $wb.Worksheets.Item("Change Log")

# Widen column B (was 5.5 characters) to fit "ENAIRE,DSNA" -> stored width 11.0
$wsLog.Columns.Item(2).ColumnWidth = 61/6

$wsLog.Range("A3").Value = 45107
$wsLog.Range("B3").Value = "ENAIRE,DSNA"
$wsLog.Range("C3").Value = 2022
$wsLog.Range("D3").Value = "55 min.  of regulation 'T21316E'  (16/09/2022) reallocated from ENAIRE to DSNA"

$wsLog.Range("A4").Value = 45107
$wsLog.Range("B4").Value = "ENAIRE,DSNA"
$wsLog.Range("C4").Value = 2022
$wsLog.Range("D4").Value = [char]0x2018 + "PAU23M" + [char]0x2019
$wsLog.Range("D4").Value = "274 min.  of regulation " + [char]0x2018 + "PAU23M" + [char]0x2019 + " (23/12/2022) reallocated from ENAIRE to DSNA"

$wsLog.Range("A5").Value = 45107
$wsLog.Range("B5").Value = "ENAIRE,DSNA"
$wsLog.Range("C5").Value = 2022
$wsLog.Range("D5").Value = "17 min. of regulation " + [char]0x2018 + "LPAU23" + [char]0x2019 + " (23/12/2022) reallocated from ENAIRE to DSNA"

Write-Host "change log updated"
